$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test zakresow")

# New headers for the two additional columns (channel multiplier, and its
# "paste values" snapshot) introduced when encoder-channel settings moved
# from the GUI loop to the acquisition loop. Match the wrap-text style
# used by the rest of the header row (B4:K4).
$ws.Range("L4").Value = "mnożnik"
$ws.Range("L4").WrapText = $true
$ws.Range("M4").Value = "mnożlik - wartość wklejona"
$ws.Range("M4").WrapText = $true

# Column L: mnoznik = I/H (row 5 entered individually, 6:12 filled as one
# block so they form a shared formula group, mirroring the source file).
$ws.Range("L5").Formula = "=I5/H5"
$ws.Range("J6:J12").Formula = "=D6*H6/I6"
$ws.Range("L6:L12").Formula = "=I6/H6"

# Column M: paste-special "values only" snapshot of column L.
$ws.Range("L5:L12").Copy()
$ws.Range("M5:M12").PasteSpecial(-4163)
$excel.CutCopyMode = $false

[void]$ws.Range("M12").Select()
